$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3957.4583
$ws.Range("I64").Value = 3441.3572
$ws.Range("K64").Value = 3441.3572
$ws.Range("M64").Value = -3193.3572
$ws.Range("H67").Value = 3957.4583
$ws.Range("I67").Value = 3441.3572
$ws.Range("K67").Value = 3441.3572
$ws.Range("M67").Value = -2583.3572
$ws.Range("H70").Value = 1126.6666
$ws.Range("I70").Value = 940
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 2820
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -2550
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 1126.6666
$ws.Range("I73").Value = 940
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 2820
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -1884
$ws.Range("N73").Value = -6372
$ws.Range("H74").Value = 12504471
$ws.Range("I74").Value = 2557
$ws.Range("J74").Value = 15629950
$ws.Range("K74").Value = 2557
$ws.Range("L74").Value = 15629950
$ws.Range("M74").Value = -1621
$ws.Range("N74").Value = -15631822
$ws.Range("H76").Value = 2139697.2
$ws.Range("I76").Value = 2972.2
$ws.Range("J76").Value = 5053413.5
$ws.Range("K76").Value = 2972.2
$ws.Range("L76").Value = 5053413.5
$ws.Range("M76").Value = -2657.2
$ws.Range("N76").Value = -5054043.5
$ws.Range("H77").Value = 12504471
$ws.Range("I77").Value = 2557
$ws.Range("J77").Value = 15629950
$ws.Range("K77").Value = 12785
$ws.Range("L77").Value = 78149750
$ws.Range("M77").Value = -8105
$ws.Range("N77").Value = -78159110
$ws.Range("H79").Value = 2139697.2
$ws.Range("I79").Value = 2972.2
$ws.Range("J79").Value = 5053413.5
$ws.Range("K79").Value = 2972.2
$ws.Range("L79").Value = 5053413.5
$ws.Range("M79").Value = -1880.2
$ws.Range("N79").Value = -5055597.5
$ws.Range("H80").Value = 4252
$ws.Range("I80").Value = 1193.375
$ws.Range("J80").Value = 6291.0835
$ws.Range("K80").Value = 3580.125
$ws.Range("L80").Value = 18873.2505
$ws.Range("M80").Value = -2582.125
$ws.Range("N80").Value = -20869.2505
$ws.Range("H82").Value = 1321
$ws.Range("I82").Value = 1321
$ws.Range("K82").Value = 3963
$ws.Range("M82").Value = -3557
$ws.Range("H83").Value = 4252
$ws.Range("I83").Value = 1193.375
$ws.Range("J83").Value = 6291.0835
$ws.Range("K83").Value = 10740.375
$ws.Range("L83").Value = 56619.7515
$ws.Range("M83").Value = -5748.375
$ws.Range("N83").Value = -66603.7515
$ws.Range("H85").Value = 1321
$ws.Range("I85").Value = 1321
$ws.Range("K85").Value = 3963
$ws.Range("M85").Value = -2559
$ws.Range("H86").Value = 8597.75
$ws.Range("I86").Value = 2199.7144
$ws.Range("J86").Value = 13574
$ws.Range("K86").Value = 2199.7144
$ws.Range("L86").Value = 13574
$ws.Range("M86").Value = -1076.7144
$ws.Range("N86").Value = -15820
$ws.Range("H89").Value = 8597.75
$ws.Range("I89").Value = 2199.7144
$ws.Range("J89").Value = 13574
$ws.Range("K89").Value = 10998.572
$ws.Range("L89").Value = 67870
$ws.Range("M89").Value = -5382.572
$ws.Range("N89").Value = -79102
$ws.Range("H94").Value = 1002
$ws.Range("I94").Value = 1002
$ws.Range("K94").Value = 1002
$ws.Range("M94").Value = -551
$ws.Range("H97").Value = 1015
$ws.Range("J97").Value = 1015
$ws.Range("L97").Value = 3045
$ws.Range("N97").Value = -4037
$ws.Range("H98").Value = 266.94736
$ws.Range("I98").Value = 266.94736
$ws.Range("K98").Value = 266.94736
$ws.Range("M98").Value = 1231.05264
$ws.Range("H100").Value = 3496.4285
$ws.Range("I100").Value = 3272.7273
$ws.Range("J100").Value = 4316.6665
$ws.Range("K100").Value = 3272.7273
$ws.Range("L100").Value = 4316.6665
$ws.Range("M100").Value = -2731.7273
$ws.Range("N100").Value = -5398.6665
$ws.Range("H122").Value = 266.94736
$ws.Range("I122").Value = 266.94736
$ws.Range("K122").Value = 800.84208
$ws.Range("M122").Value = 1649.15792
$ws.Range("H138").Value = 1630.3889
$ws.Range("I138").Value = 621.125
$ws.Range("J138").Value = 2437.8
$ws.Range("K138").Value = 1863.375
$ws.Range("L138").Value = 7313.400000000001
$ws.Range("M138").Value = 3276.625
$ws.Range("N138").Value = -17593.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2264.2896
$ws.Range("I94").Value = 1068.1
$ws.Range("J94").Value = 6750
$ws.Range("K94").Value = 1068.1
$ws.Range("L94").Value = 6750
$ws.Range("M94").Value = -617.0999999999999
$ws.Range("N94").Value = -7652
$ws.Range("H99").Value = 2060.6
$ws.Range("I99").Value = 1714.1428
$ws.Range("K99").Value = 1714.1428
$ws.Range("M99").Value = -216.1428000000001
$ws.Range("H105").Value = 3848507.5
$ws.Range("I105").Value = 1971
$ws.Range("J105").Value = 8336133.5
$ws.Range("K105").Value = 1971
$ws.Range("L105").Value = 8336133.5
$ws.Range("M105").Value = -224
$ws.Range("N105").Value = -8339627.5
$ws.Range("H134").Value = 30552.432
$ws.Range("I134").Value = 46435.25
$ws.Range("J134").Value = 1230.3077
$ws.Range("K134").Value = 139305.75
$ws.Range("L134").Value = 3690.9231
$ws.Range("M134").Value = -136770.75
$ws.Range("N134").Value = -8760.9231

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11908228
$ws.Range("I99").Value = 2856.923
$ws.Range("K99").Value = 2856.923
$ws.Range("M99").Value = -1358.923
$ws.Range("H107").Value = 1213.625
$ws.Range("I107").Value = 590.25
$ws.Range("J107").Value = 1837
$ws.Range("K107").Value = 590.25
$ws.Range("L107").Value = 1837
$ws.Range("M107").Value = 1329.75
$ws.Range("N107").Value = -5677
$ws.Range("H126").Value = 11908228
$ws.Range("I126").Value = 2856.923
$ws.Range("K126").Value = 8570.769
$ws.Range("M126").Value = -6100.769
$ws.Range("H132").Value = 19481.268
$ws.Range("I132").Value = 21897.04
$ws.Range("J132").Value = 7402.4
$ws.Range("K132").Value = 65691.12
$ws.Range("L132").Value = 22207.2
$ws.Range("M132").Value = -63161.12
$ws.Range("N132").Value = -27267.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9858.091
$ws.Range("I113").Value = 50450.5
$ws.Range("J113").Value = 837.55554
$ws.Range("K113").Value = 151351.5
$ws.Range("L113").Value = 2512.66662
$ws.Range("M113").Value = -149181.5
$ws.Range("N113").Value = -6852.66662
$ws.Range("H131").Value = 747.8182
$ws.Range("J131").Value = 754.6
$ws.Range("L131").Value = 2263.8
$ws.Range("N131").Value = -12343.8
$ws.Range("H132").Value = 1067.1177
$ws.Range("J132").Value = 1260.2222
$ws.Range("L132").Value = 11341.9998
$ws.Range("N132").Value = -16401.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 12505600
$ws.Range("I52").Value = 2800
$ws.Range("J52").Value = 14291714
$ws.Range("K52").Value = 2800
$ws.Range("L52").Value = 14291714
$ws.Range("M52").Value = -2541
$ws.Range("N52").Value = -14292232
$ws.Range("H132").Value = 61562.96
$ws.Range("I132").Value = 51651.855
$ws.Range("J132").Value = 103189.6
$ws.Range("K132").Value = 154955.565
$ws.Range("L132").Value = 309568.8
$ws.Range("M132").Value = -152425.565
$ws.Range("N132").Value = -314628.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 26720.5
$ws.Range("I34").Value = 10324.75
$ws.Range("K34").Value = 10324.75
$ws.Range("M34").Value = -10152.75
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H82").Value = 1683.3704
$ws.Range("I82").Value = 2081.5833
$ws.Range("J82").Value = 1364.8
$ws.Range("K82").Value = 2081.5833
$ws.Range("L82").Value = 1364.8
$ws.Range("M82").Value = -1720.5833
$ws.Range("N82").Value = -2086.8
$ws.Range("H85").Value = 1683.3704
$ws.Range("I85").Value = 2081.5833
$ws.Range("J85").Value = 1364.8
$ws.Range("K85").Value = 2081.5833
$ws.Range("L85").Value = 1364.8
$ws.Range("M85").Value = -833.5832999999998
$ws.Range("N85").Value = -3860.8
$ws.Range("H93").Value = 1984.5
$ws.Range("I93").Value = 1984.5
$ws.Range("K93").Value = 1984.5
$ws.Range("M93").Value = -736.5
$ws.Range("H100").Value = 3068.0908
$ws.Range("I100").Value = 2185.7144
$ws.Range("K100").Value = 2185.7144
$ws.Range("M100").Value = -1644.7144
$ws.Range("H132").Value = 1951.08
$ws.Range("I132").Value = 1291.3572
$ws.Range("K132").Value = 3874.0716
$ws.Range("M132").Value = -1344.0716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 999
$ws.Range("I126").Value = 999
$ws.Range("J126").Value = 999
$ws.Range("K126").Value = 2997
$ws.Range("L126").Value = 2997
$ws.Range("M126").Value = -527
$ws.Range("N126").Value = -7937
